# Edit script: Horarios Linea 141 - update scraped schedule data (run 11:15:53)
$wb = $excel.ActiveWorkbook

# --- Sheet 1: LP1912 ---
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 11:15:53"
$ws1.Range("A3").Value = "Total filas: 204"

$sheet1Rows = @(
    @(118, '08:33:47', '09:34', '16_SANTA ANA', 61, 'LP1912'),
    @(119, '08:54:42', '09:34', '23_HERNANDEZ', 40, 'LP1912'),
    @(167, '11:15:53', '11:15', '16_SANTA ANA', 0, 'LP1912'),
    @(168, '11:15:53', '11:15', '86_EST CHICA-ESC AGRARIA', 0, 'LP1912'),
    @(169, '11:15:53', '11:16', '15_ABASTO', 1, 'LP1912'),
    @(170, '09:25:30', '11:19', '86_EST CHICA-ESC AGRARIA', 114, 'LP1912'),
    @(171, '11:15:53', '11:20', '26_HERNANDEZ', 5, 'LP1912'),
    @(172, '10:11:11', '11:20', '86_EST CHICA-ESC AGRARIA', 69, 'LP1912'),
    @(173, '09:25:30', '11:21', '26_HERNANDEZ', 116, 'LP1912'),
    @(174, '11:15:53', '11:26', '225_C ROCA-H SUR', 11, 'LP1912'),
    @(175, '10:11:11', '11:27', '225_C ROCA-H SUR', 76, 'LP1912'),
    @(176, '10:11:11', '11:32', '81_EL PELIGRO', 81, 'LP1912'),
    @(177, '10:50:37', '11:34', '23_HERNANDEZ', 44, 'LP1912'),
    @(178, '10:50:37', '11:35', '11_ETCHEVERRY', 45, 'LP1912'),
    @(179, '11:15:53', '11:35', '23_HERNANDEZ', 20, 'LP1912'),
    @(180, '10:11:11', '11:38', '10_OLMOS', 87, 'LP1912'),
    @(181, '10:50:37', '11:41', '17_ROMERO', 51, 'LP1912'),
    @(182, '10:11:11', '11:42', '17_ROMERO', 91, 'LP1912'),
    @(183, '10:50:37', '11:43', '10_OLMOS', 53, 'LP1912'),
    @(184, '10:11:11', '11:51', '215B_EL PATO', 100, 'LP1912'),
    @(185, '11:15:53', '11:58', '225_GOMEZ', 43, 'LP1912'),
    @(186, '10:11:11', '11:59', '225_GOMEZ', 108, 'LP1912'),
    @(187, '10:11:11', '12:02', '84_COLONIA URQUIZA-ESC 49', 111, 'LP1912'),
    @(188, '11:15:53', '12:04', '23_HERNANDEZ', 49, 'LP1912'),
    @(189, '10:50:37', '12:06', '16_P MOR-SANTA ANA', 76, 'LP1912'),
    @(190, '10:50:37', '12:06', '14_ABASTO', 76, 'LP1912'),
    @(191, '10:50:37', '12:07', '10_OLMOS', 77, 'LP1912'),
    @(192, '10:11:11', '12:07', '14_ABASTO', 116, 'LP1912'),
    @(193, '10:11:11', '12:07', '16_P MOR-SANTA ANA', 116, 'LP1912'),
    @(194, '11:15:53', '12:12', '10_OLMOS', 57, 'LP1912'),
    @(195, '10:50:37', '12:20', '215A_EL PATO', 90, 'LP1912'),
    @(196, '11:15:53', '12:20', '26_HERNANDEZ', 65, 'LP1912'),
    @(197, '11:15:53', '12:20', '14_ABASTO', 65, 'LP1912'),
    @(198, '10:50:37', '12:21', '26_HERNANDEZ', 91, 'LP1912'),
    @(199, '10:50:37', '12:21', '14_ABASTO', 91, 'LP1912'),
    @(200, '10:50:37', '12:22', '17_ROMERO', 92, 'LP1912'),
    @(201, '11:15:53', '12:34', '11_ETCHEVERRY', 79, 'LP1912'),
    @(202, '10:50:37', '12:36', '27_EL RETIRO', 106, 'LP1912'),
    @(203, '10:50:37', '12:38', '17_179 Y 38', 108, 'LP1912'),
    @(204, '11:15:53', '12:40', '10_OLMOS', 85, 'LP1912'),
    @(205, '11:15:53', '12:46', '17_ROMERO', 91, 'LP1912'),
    @(206, '11:15:53', '12:48', '11_ETCHEVERRY', 93, 'LP1912'),
    @(207, '11:15:53', '13:02', '15_ABASTO', 107, 'LP1912'),
    @(208, '11:15:53', '13:06', '16_P MOR-SANTA ANA', 111, 'LP1912'),
    @(209, '11:15:53', '13:13', '215D_EL PATO', 118, 'LP1912')
)

foreach ($row in $sheet1Rows) {
    $r = $row[0]
    $ws1.Range("A$r").Value = $row[1]
    $ws1.Range("B$r").Value = $row[2]
    $ws1.Range("C$r").Value = $row[3]
    $ws1.Range("D$r").Value = $row[4]
    $ws1.Range("E$r").Value = $row[5]
}

# --- Sheet 2: LP1912-215 ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 11:15:53"
$ws2.Range("A3").Value = "Total filas: 25"

$ws2.Range("A30").Value = "11:15:53"
$ws2.Range("B30").Value = "13:13"
$ws2.Range("C30").Value = "215D_EL PATO"
$ws2.Range("D30").Value = 118
$ws2.Range("E30").Value = "LP1912"

# --- Sheet 3: 6203-6173 ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 11:15:53"
$ws3.Range("A3").Value = "Total filas: 32"

# Existing row 35 (10:11:11 / 12:04 / 215A_LA PLATA / 113 / L6173) shifts down to row 36.
$ws3.Range("A36").Value = "10:11:11"
$ws3.Range("B36").Value = "12:04"
$ws3.Range("C36").Value = "215A_LA PLATA"
$ws3.Range("D36").Value = 113
$ws3.Range("E36").Value = "L6173"

# New row 35 inserted before it.
$ws3.Range("A35").Value = "11:15:53"
$ws3.Range("B35").Value = "11:15"
$ws3.Range("C35").Value = "215C_LA PLATA"
$ws3.Range("D35").Value = 0
$ws3.Range("E35").Value = "L6203"

# New row 37 appended after.
$ws3.Range("A37").Value = "11:15:53"
$ws3.Range("B37").Value = "12:53"
$ws3.Range("C37").Value = "215C_LA PLATA"
$ws3.Range("D37").Value = 98
$ws3.Range("E37").Value = "L6203"
